# Se crea la hoja "Agendas" (respaldo de la clase AgendaSemanal) para
# guardar las agendas semanales de los asesores en el archivo de Excel.

$wb = $excel.ActiveWorkbook

# La nueva hoja se ubica entre "asesorias" y "notificaciones".
$wsAsesorias = $wb.Worksheets.Item("asesorias")
$wsAgendas = $wb.Worksheets.Add($null, $wsAsesorias)
$wsAgendas.Name = "Agendas"

# Encabezados de la agenda semanal: usuario del asesor + los 7 días.
$headers = @("Usuario Asesor", "Lunes", "Martes", "Miércoles ", "Jueves", "Viernes", "Sábado", "Domingo")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $wsAgendas.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
}

$wsAgendas.PageSetup.Orientation = 1

# La hoja recien creada queda como la hoja activa del libro.
$wsAgendas.Activate()
